# TMTC0032668 - Changed Test Data for LV Activities - 16 Sep 2024
#
# The "Contact" test-data row is updated to use an activity-specific
# external contact / company, and the workbook is left with the
# "Contact" sheet active (instead of "Followup"), matching the selection
# that was in effect when the author made the edit.

$wb = $excel.ActiveWorkbook

$contact = $wb.Worksheets.Item("Contact")
$contact.Range("A2").Value = "Activity Test External Contact"
$contact.Range("B2").Value = "ActivityCompany"

# Make Contact the active sheet / tab, with A2:B2 selected (active cell A2)
$contact.Activate()
$contact.Range("A2:B2").Select()
